$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Tompkins health update - today's row (day 11) added.
# Quarantine numbers (In_quarantine / Released_quarantine) are no longer
# reported -> #N/A. Other data for today (Pending/Positive/Negative/
# Total_tested/Deaths) is filled in; formatting is carried down from the
# previous "today" row (row 11) / the regular data rows above it.

# Day number
$ws.Range("A12").Value = 11

# Date - carry down the date style used on row 11 (the previous "today" row)
$ws.Range("B11").Copy()
$ws.Range("B12").PasteSpecial(-4122)
$ws.Range("B12").Formula = "=B11+1"

# Pending / Positive / Negative / Total_tested - carry down the plain
# wrap/vertical-center formatting used across the data rows
$ws.Range("C2:F2").Copy()
$ws.Range("C12:F12").PasteSpecial(-4122)

$ws.Range("C12").Value = 628
$ws.Range("D12").Value = 48
$ws.Range("E12").Value = 515
$ws.Range("F12").Value = 1191

# In_quarantine / Released_quarantine - no longer reported
$ws.Range("G12").Value = "#N/A"
$ws.Range("H12").Value = "#N/A"

# Deaths
$ws.Range("I2").Copy()
$ws.Range("I12").PasteSpecial(-4122)
$ws.Range("I12").Value = 0

$ws.Application.CutCopyMode = $false

# Move the active selection down to the next blank row, as happens after
# entering a new row of data
$ws.Range("H13").Select()

$wb.Save()
